$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 56, shifting existing rows 56:157 down to 57:158
# (this also grows the used range to A1:R158, matching the dimension change).
$ws.Rows(56).Insert()

# Populate the newly inserted row 56 with this week's new record.
$ws.Range("A56").Value = 10
$ws.Range("B56").Value = 'Vega Modelo de Temuco'
$ws.Range("C56").Value = 'La Araucanía'
$ws.Range("D56").Value = 44519
$ws.Range("E56").Value = 9
$ws.Range("F56").Value = 100112005
$ws.Range("G56").Value = 'Puerro'
$ws.Range("H56").Value = 'Azul de Maquehue'
$ws.Range("I56").Value = 'Primera'
$ws.Range("J56").Value = 155
$ws.Range("K56").Value = 7000
$ws.Range("L56").Value = 7000
$ws.Range("M56").Value = 7000
$ws.Range("N56").Value = '$/docena de paquetes'
$ws.Range("O56").Value = 'Provincia de Cautín'
$ws.Range("P56").Value = 583
$ws.Range("Q56").Value = 12
$ws.Range("R56").Value = 'Hortaliza'
